$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The sheet currently holds 23 data rows (1..23). The edit inserts a new
# "Docentes responsaveis:" block (1 label row + 4 professor rows = 5 rows)
# right after the "Objectives:" row (old row 11), pushing the old rows
# 12..23 down to 17..28. It also corrects several mis-aligned text values
# (several B/C cells previously held the wrong text, e.g. a professor's
# name where the syllabus text belongs) and replaces the Objetivos (PT)
# text with its real content.
# -----------------------------------------------------------------------

# 1) Insert 5 blank rows at row 12 (shifts old rows 12-23 -> 17-28)
$ws.Rows.Item(12).Resize(5).Insert()

# Newly inserted rows 12-16 copy the formatting (and empty styled cells)
# of the row above for columns A,B,C. Remove the cells that must not
# exist at all in the final layout (row 12 only has col A; rows 13-16
# only have cols B and C).
$ws.Range("B12:C12").Clear()
$ws.Range("A13:A16").Clear()

# 2) Fix row 10 ("Objetivos:") text - it had a misplaced professor name
$objetivosPt = 'Estudo formal da teoria dos campos eletromagnéticos independentes do tempo ou para situações quase-estáticas. Teoria das ondas eletromagnéticas.'
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# Row 11 ("Objectives:") text is unchanged (Formal study...) - leave as-is

# 3) Fill the new "Docentes responsaveis:" block (rows 12-16)
$ws.Range("A12").Value = 'Docentes responsáveis:'

$prof1 = '5840726 - Cristina Bormio Nunes'
$ws.Range("B13").Value = $prof1
$ws.Range("C13").Value = $prof1

$prof2 = '6495737 - Durval Rodrigues Junior'
$ws.Range("B14").Value = $prof2
$ws.Range("C14").Value = $prof2

$prof3 = '1341653 - Maria José Ramos Sandim'
$ws.Range("B15").Value = $prof3
$ws.Range("C15").Value = $prof3

$prof4 = '1643715 - Paulo Atsushi Suzuki'
$ws.Range("B16").Value = $prof4
$ws.Range("C16").Value = $prof4

# 4) Correct the text of the shifted rows (now 17-25) that previously held
#    misplaced content (professor names sitting in the wrong cells).

# Row 17 (was row 12): "Programa resumido:"
$programaResumidoPt = 'Eletrostática. Magnetostática. Campos variantes no tempo. Equações de Maxwell. Ondas eletromagnéticas.'
$ws.Range("B17").Value = $programaResumidoPt
$ws.Range("C17").Value = $programaResumidoPt

# Row 18 (was row 13): "Short syllabus:" - text unchanged, already correct

# Row 19 (was row 14): "Programa:"
$programaPt = 'Eletrostática (campo eletrostático; potencial elétrico; trabalho e energia em eletrostática).  Técnicas especiais para a resolução da equação de Laplace (método das imagens; separação de variáveis). Campo elétrico da matéria (polarização elétrica; campo de objeto polarizado; cargas ligadas; deslocamento elétrico; dielétricos (lineares). Magnetostática (Lei de Lorentz; Lei de Biot-Savart; Lei de Ampére; vetor potencial magnético).  Campo magnético na matéria (magnetização; campos de objeto magnetizado; campo auxiliar H; Eletrodinâmica (força eletromotriz; indução eletromagnética; equações de Maxwell; lei de conservação de carga). Ondas eletromagnéticas (propagação no vácuo e na matéria; reflexão e transmissão), equação de ondas (planas)  e condições de contorno (interfaces). Radiação de dipolo elétrico.'
$ws.Range("B19").Value = $programaPt
$ws.Range("C19").Value = $programaPt

# Row 20 (was row 15): "Syllabus:" - text unchanged, already correct

# Row 21 (was row 16): "Avaliação:" label only - unchanged

# Row 22 (was row 17): "Método:"
$metodoPt = 'Aulas expositivas e  exercícios comentados'
$ws.Range("B22").Value = $metodoPt
$ws.Range("C22").Value = $metodoPt

# Row 23 (was row 18): "Critério:"
$criterioPt = 'Média final calculada pelas notas de 2 provas (P1 e P2), seguindo os pesos MF=(P1+2*P2)/3, ou seja, peso 1 para a P1 e peso 2 para a P2.'
$ws.Range("B23").Value = $criterioPt
$ws.Range("C23").Value = $criterioPt

# Row 24 (was row 19): "Norma de recuperação:"
$normaPt = 'Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("B24").Value = $normaPt
$ws.Range("C24").Value = $normaPt

# Row 25 (was row 20): "Bibliografia:"
$bibliografiaPt = "CHENG,DAVID K.Field and Wave Electromagnetics. Addison Weslwy Publishing Company.1989.`nSLATER, J.C.; FRANK, N.H. Electromagnetism. McGraw-Hill, New York, 1974.`nMARION, J.B. Classical Electromagnetic Radiation. Academic Press, New York, 1965.`nBOHN, E.V. Introduction to electromagnetic fields and waves. Addison Wesley, 1968. `nREITZ, J.R.; MILFORD, F.J. Foundations of eletromagnetic theory. Addison Wesley, Publishing, Co. 1970. GRIFFITHS, D.J. Introduction to Electrodynamics. Prentice Hall, New York. 1998. `nRAMO, WHINNERY E VAN DUZER, Fields and Waves in Communication Electronics, Wiley."
$ws.Range("B25").Value = $bibliografiaPt
$ws.Range("C25").Value = $bibliografiaPt

# Row 26 (was row 21): "Requisitos:" label only - unchanged
# Rows 27-28 (was rows 22-23): LOB1052/LOB1053 requirement text - unchanged

Write-Host "Done applying LOM3205 updates."
